$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览": remove the 2 events that were dropped from the source feed
#   - row 2 (A=1): 昆山·M·M国风动漫展
#   - row 4 (A=3, originally): 张家港·万达国庆动漫随宅自由展
# Delete bottom-most row first so the earlier row index stays valid.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Rows.Item(4).Delete()
$wsExpo.Rows.Item(2).Delete()

# Renumber the leading index column (A) back to a clean 1..N sequence.
for ($i = 1; $i -le 21; $i++) {
    $wsExpo.Cells.Item($i + 1, 1).Value = $i
}

# Refresh "想去人数" (col F) counters that moved since the last scrape.
$expoFUpdates = @{
    1  = 749
    2  = 662
    4  = 92
    5  = 1177
    6  = 20
    7  = 30
    8  = 41
    9  = 574
    15 = 81
    16 = 287
    17 = 392
    18 = 485
    19 = 124
    20 = 5917
    21 = 5285
}
foreach ($row in $expoFUpdates.Keys) {
    $wsExpo.Cells.Item($row + 1, 6).Value = $expoFUpdates[$row]
}

# ---------------------------------------------------------------------------
# Sheet "演出": bump the same counter refresh for the one row it carries.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(2, 6).Value = 88

# ---------------------------------------------------------------------------
# Sheet "全部类型": same drop as "展览" (same 2 removed events) plus the same
# counter refreshes (including the 演出 row embedded in this combined sheet).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows.Item(4).Delete()
$wsAll.Rows.Item(2).Delete()

for ($i = 1; $i -le 24; $i++) {
    $wsAll.Cells.Item($i + 1, 1).Value = $i
}

$allFUpdates = @{
    1  = 749
    2  = 662
    4  = 92
    5  = 1177
    6  = 20
    7  = 30
    8  = 41
    9  = 574
    12 = 88
    17 = 81
    18 = 287
    19 = 392
    20 = 485
    21 = 124
    22 = 5917
    24 = 5285
}
foreach ($row in $allFUpdates.Keys) {
    $wsAll.Cells.Item($row + 1, 6).Value = $allFUpdates[$row]
}
